$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 77, shifting rows 77-86 down to 80-89.
$ws.Range("A77:A79").EntireRow.Insert()

# Row 77: shared columns
$ws.Range("A77").Value = 1
$ws.Range("B77").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C77").Value = "Arica y Parinacota"
$ws.Range("D77").Value = 45106
$ws.Range("E77").Value = 15
$ws.Range("F77").Value = 100112045
$ws.Range("G77").Value = "Zapallo"
$ws.Range("H77").Value = "Camote"
$ws.Range("I77").Value = "1a (guarda)"
$ws.Range("J77").Value = 400
$ws.Range("K77").Value = 430
$ws.Range("L77").Value = 450
$ws.Range("M77").Value = 435
$ws.Range("N77").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O77").Value = "Región de O'Higgins"
$ws.Range("P77").Value = 435
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"

# Row 78
$ws.Range("A78").Value = 1
$ws.Range("B78").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C78").Value = "Arica y Parinacota"
$ws.Range("D78").Value = 45106
$ws.Range("E78").Value = 15
$ws.Range("F78").Value = 100112045
$ws.Range("G78").Value = "Zapallo"
$ws.Range("H78").Value = "Camote"
$ws.Range("I78").Value = "2a (guarda)"
$ws.Range("J78").Value = 250
$ws.Range("K78").Value = 400
$ws.Range("L78").Value = 430
$ws.Range("M78").Value = 418
$ws.Range("N78").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O78").Value = "Región de O'Higgins"
$ws.Range("P78").Value = 418
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"

# Row 79
$ws.Range("A79").Value = 1
$ws.Range("B79").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C79").Value = "Arica y Parinacota"
$ws.Range("D79").Value = 45106
$ws.Range("E79").Value = 15
$ws.Range("F79").Value = 100112045
$ws.Range("G79").Value = "Zapallo"
$ws.Range("H79").Value = "Camote"
$ws.Range("I79").Value = "3a (guarda)"
$ws.Range("J79").Value = 250
$ws.Range("K79").Value = 350
$ws.Range("L79").Value = 400
$ws.Range("M79").Value = 380
$ws.Range("N79").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O79").Value = "Región de O'Higgins"
$ws.Range("P79").Value = 380
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"

Write-Host "Rows 77-79 inserted and populated; dimension now extends to row 89."
